# Maj docs suivi de projet - suivi perso / suivi taches
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 31 - new task entry
$ws.Range("B31").Value = "Travail sur les textures photoshop"
$ws.Range("C31").NumberFormat = "d-mmm"
$ws.Range("C31").Value = 42329
$ws.Range("D31").Value = 4

# Row 32 - new task entry
$ws.Range("B32").Value = "Reprise fichiers templates pour responsive design"
$ws.Range("C32").NumberFormat = "d-mmm"
$ws.Range("C32").Value = 42331
$ws.Range("D32").Value = 1

# Row 33 - same task, different day/hours
$ws.Range("B33").Value = "Reprise fichiers templates pour responsive design"
$ws.Range("C33").NumberFormat = "d-mmm"
$ws.Range("C33").Value = 42336
$ws.Range("D33").Value = 5

# Row 34 - new task entry
$ws.Range("B34").Value = "Modification organisation arborsence application"
$ws.Range("C34").NumberFormat = "d-mmm"
$ws.Range("C34").Value = 42336
$ws.Range("D34").Value = 0.5

# Update the active selection on the sheet
$ws.Range("H24").Select()
